# Fix cost bug: update calibrated production-initialization values for
# several IPPU subsectors on rows 96-104 and 111-112. Each affected row
# holds a single constant value repeated across columns J:AS (the 36
# trajectory-year columns), so we overwrite that whole range per row with
# the corrected constant.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    96  = 3944281.371
    97  = 391247.6029
    98  = 193163.7998
    99  = 29174.73605
    100 = 4136.366244
    101 = 1660792.974
    103 = 803025.1581999999
    104 = 386455.3512
    111 = 13449.72472
    112 = 55827.44077
}

foreach ($row in $updates.Keys) {
    $value = $updates[$row]
    $range = $ws.Range("J$row`:AS$row")
    $range.Value = $value
}
